# Update "want-to-go" counts (column F) on the "展览" (Exhibition) and
# "全部类型" (All Types) sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition listing) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 286   # was 285
$ws1.Range("F3").Value = 96    # was 95
$ws1.Range("F4").Value = 1101  # was 1091
$ws1.Range("F5").Value = 576   # was 571

# --- Sheet "全部类型" (combined listing of all event types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 286   # was 285
$ws4.Range("F3").Value = 96    # was 95
$ws4.Range("F4").Value = 1101  # was 1091
$ws4.Range("F6").Value = 576   # was 571
